# Update cryptocurrency price and volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format first, otherwise Excel auto-converts them to
# floating point numbers instead of keeping the original text formatting.
$numericLookingCells = @("D5", "D6", "D7", "D9", "D10", "D14", "D15", "D17", "D18", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D32", "D33", "D35", "D36", "D38", "D40", "D41", "D42", "D45", "D46", "D47", "D51")
foreach ($cellAddr in $numericLookingCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = '51.564.73'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '3.047.13'
$ws.Range("E3").Value = '  +2.91%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '385.28'
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").Value = '102.75'
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("D7").Value = '0.544'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.584'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("D10").Value = '36.83'
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '3.530.55'
$ws.Range("D14").Value = '18.65'
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").Value = '7.75'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").Value = '3.073.40'
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("D17").Value = '0.974'
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").Value = '10.53'
$ws.Range("E18").Value = '  -5.61%  '
$ws.Range("D19").Value = '51.640.94'
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").Value = '3.15'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").Value = '12.41'
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").Value = '0.0₃0963'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").Value = '70.24'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = '268.46'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").Value = '3.16'
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("D26").Value = '8.24'
$ws.Range("E26").Value = '  +5.71%  '
$ws.Range("D27").Value = '27.07'
$ws.Range("E27").Value = '  +4.53%  '
$ws.Range("E28").Value = '  +3.95%  '
$ws.Range("D29").Value = '7.25'
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("E31").Value = '  -1.67%  '
$ws.Range("D32").Value = '10.27'
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '34.71'
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("E34").Value = '  +2.83%  '
$ws.Range("D35").Value = '50.45'
$ws.Range("E35").Value = '  -1.41%  '
$ws.Range("D36").Value = '0.0445'
$ws.Range("E36").Value = '  +1.81%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = '3.35'
$ws.Range("E38").Value = '  +2.89%  '
$ws.Range("E39").Value = '  +7.56%  '
$ws.Range("D40").Value = '16.98'
$ws.Range("E40").Value = '  +2.80%  '
$ws.Range("D41").Value = '1.87'
$ws.Range("E41").Value = '  +2.08%  '
$ws.Range("D42").Value = '128.12'
$ws.Range("E42").Value = '  +2.33%  '
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").Value = '3.75'
$ws.Range("E45").Value = '  +5.22%  '
$ws.Range("D46").Value = '21.90'
$ws.Range("E46").Value = '  +1.88%  '
$ws.Range("D47").Value = '2.48'
$ws.Range("E47").Value = '  +4.59%  '
$ws.Range("E48").Value = '  +2.86%  '
$ws.Range("D49").Value = '2.035.26'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").Value = '3.347.09'
$ws.Range("E50").Value = '  +2.89%  '
$ws.Range("D51").Value = '0.205'
$ws.Range("E51").Value = '  +6.27%  '

# Restore default (General) style on the cells we forced to Text so the
# cell formatting matches the rest of the sheet (no explicit style index).
foreach ($cellAddr in $numericLookingCells) {
    $ws.Range($cellAddr).Style = "Normal"
}
